# Apply "Add descriptions titles" commit:
#  - Metadata sheet: fill in Title / Date / Description values
#  - Elements sheet: fill in the root Extension row's Short / Definition,
#    and clear the (now redundant) "N/A" RIM mapping for that row.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B5").Value  = "DMI IP Id logiciel"
$meta.Range("B8").Value  = "2026-02-25T08:15:31+00:00"
$meta.Range("B12").Value = "Extension créée dans ce volet pour représenter l'IP Id logiciel."

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("L2").Value  = "DMI IP Id logiciel"
$elements.Range("M2").Value  = "Extension créée dans ce volet pour représenter l'IP Id logiciel."
$elements.Range("AK2").Value = ""
